# Week 15 simulations update
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Rushing")
$ws2 = $wb.Worksheets.Item("Receiving")

# ---------- Rushing sheet updates ----------
# Row 4 (D.Henderson)
$ws1.Range("C4").Value = 0
$ws1.Range("D4").Value = 0
$ws1.Range("E4").Value = 0
$ws1.Range("F4").Value = 0

# Row 5 (S.Michel)
$ws1.Range("C5").Value = 64
$ws1.Range("D5").Value = 30
$ws1.Range("E5").Value = 12
$ws1.Range("F5").Value = 25

# Row 8 (M.Sargent)
$ws1.Range("C8").Value = 10
$ws1.Range("D8").Value = 10
$ws1.Range("E8").Value = 5
$ws1.Range("F8").Value = 7

# New row 10 (V.Jefferson) - copy formatting from row 9's A cell first
[void]$ws1.Range("A9").Copy()
[void]$ws1.Range("A10").PasteSpecial(-4122)
$ws1.Range("A10").Value = 8
$ws1.Range("B10").Value = "V.Jefferson"
$ws1.Range("C10").Value = 1
$ws1.Range("D10").Value = 0
$ws1.Range("E10").Value = 0
$ws1.Range("F10").Value = 0

# ---------- Receiving sheet updates ----------
# Row 2 (D.Henderson)
$ws2.Range("C2").Value = 0
$ws2.Range("D2").Value = 0
$ws2.Range("E2").Value = 0
$ws2.Range("F2").Value = 0
$ws2.Range("G2").Value = 0
$ws2.Range("H2").Value = 0

# Row 3 (S.Michel)
$ws2.Range("C3").Value = 18

# Row 4 (C.Kupp)
$ws2.Range("C4").Value = 116
$ws2.Range("D4").Value = 93
$ws2.Range("E4").Value = 38
$ws2.Range("F4").Value = 20
$ws2.Range("G4").Value = 30
$ws2.Range("H4").Value = 21

# Row 5 (V.Jefferson)
$ws2.Range("C5").Value = 52
$ws2.Range("D5").Value = 31
$ws2.Range("E5").Value = 23
$ws2.Range("F5").Value = 10
$ws2.Range("G5").Value = 14

# Row 6 (B.Skowronek)
$ws2.Range("C6").Value = 18

# Row 7 (O.Beckham)
$ws2.Range("C7").Value = 0
$ws2.Range("D7").Value = 0
$ws2.Range("E7").Value = 0
$ws2.Range("F7").Value = 0
$ws2.Range("G7").Value = 0
$ws2.Range("H7").Value = 0

# Row 8 (K.Blanton)
$ws2.Range("C8").Value = 4
$ws2.Range("D8").Value = 2
$ws2.Range("E8").Value = 1
$ws2.Range("F8").Value = 1

# ---------- Selection / active sheet ----------
[void]$ws1.Range("G8").Select()
[void]$ws2.Activate()
[void]$ws2.Range("H3").Select()
